$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.12
$ws.Range("G2").Value = 2.16
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.55
$ws.Range("N2").Value = 5.8
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 2.64
$ws.Range("Q2").Value = 1.55
$ws.Range("R2").Value = 1.68
$ws.Range("S2").Value = 2.32
$ws.Range("T2").Value = 1.54
$ws.Range("U2").Value = 2.72
$ws.Range("V2").Value = 1.4
$ws.Range("W2").Value = 1.86
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 21
$ws.Range("AB2").Value = 15.5
$ws.Range("AC2").Value = 10
$ws.Range("AN2").Value = 9.800000000000001

# Row 3
$ws.Range("F3").Value = 2.02
$ws.Range("G3").Value = 2.12
$ws.Range("N3").Value = 4.9
$ws.Range("P3").Value = 2.34
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 1.54
$ws.Range("S3").Value = 2.6
$ws.Range("T3").Value = 1.64
$ws.Range("U3").Value = 2.44
$ws.Range("W3").Value = 1.89
$ws.Range("Y3").Value = 23
$ws.Range("AB3").Value = 15
$ws.Range("AF3").Value = 17.5
$ws.Range("AK3").Value = 24
$ws.Range("AL3").Value = 34
$ws.Range("AN3").Value = 12.5

# Row 4
$ws.Range("F4").Value = 5.5
$ws.Range("I4").Value = 1.73
$ws.Range("K4").Value = 4.4
$ws.Range("L4").Value = 1.31
$ws.Range("N4").Value = 3.9
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.38
$ws.Range("S4").Value = 3.15
$ws.Range("V4").Value = 2.36
$ws.Range("AA4").Value = 18
$ws.Range("AI4").Value = 38
$ws.Range("AO4").Value = 10.5

# Row 6
$ws.Range("H6").Value = 2.7
$ws.Range("I6").Value = 2.84
$ws.Range("K6").Value = 3.2
$ws.Range("N6").Value = 2.48
$ws.Range("T6").Value = 2.2
$ws.Range("AB6").Value = 8.4
$ws.Range("AD6").Value = 16

# Row 8
$ws.Range("T8").Value = 2.28
$ws.Range("X8").Value = 14.5
$ws.Range("AA8").Value = 570
$ws.Range("AH8").Value = 36
$ws.Range("AO8").Value = 470

# Row 9
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 2.82
$ws.Range("J9").Value = 3.15
$ws.Range("M9").Value = 1.12
$ws.Range("N9").Value = 2.7
$ws.Range("O9").Value = 1.55
$ws.Range("Q9").Value = 2.68
$ws.Range("R9").Value = 1.2
$ws.Range("S9").Value = 5.6
$ws.Range("T9").Value = 2.14
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.54
$ws.Range("W9").Value = 1.5
$ws.Range("X9").Value = 8.800000000000001
$ws.Range("Y9").Value = 8.4
$ws.Range("Z9").Value = 16
$ws.Range("AA9").Value = 46
$ws.Range("AB9").Value = 8.800000000000001
$ws.Range("AC9").Value = 7
$ws.Range("AF9").Value = 17.5
$ws.Range("AK9").Value = 46

# Row 10
$ws.Range("F10").Value = 1.37
$ws.Range("G10").Value = 1.41
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 5.3
$ws.Range("N10").Value = 3.45
$ws.Range("O10").Value = 1.37
$ws.Range("P10").Value = 1.82
$ws.Range("S10").Value = 3.85
$ws.Range("V10").Value = 1.08
$ws.Range("W10").Value = 3.4
$ws.Range("Y10").Value = 30
$ws.Range("AA10").Value = 810
$ws.Range("AE10").Value = 340
$ws.Range("AL10").Value = 65
$ws.Range("AO10").Value = 720

# Row 11
$ws.Range("M11").Value = 1.05
